# Actualización automática de catálogo y fotos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("catalogo")

# Row 28: shift image filenames left (E empty -> imagen1, F imagen1 -> imagen2, G imagen2 -> cleared)
$ws.Range("E28").Value = "sf beige puesto 1.png"
$ws.Range("F28").Value = "sf beige puesto 2.png"
$ws.Range("G28").Value = ""

# Row 34: shift image filenames left (E empty -> imagen, F imagen -> cleared)
$ws.Range("E34").Value = "pontevedra verde 2.png"
$ws.Range("F34").Value = ""

# Row 38: shift image filenames left (F empty -> imagen, G imagen -> cleared)
$ws.Range("F38").Value = "mallorca cuerda verde 3.JPG"
$ws.Range("G38").Value = ""
